$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (and its "through" date label) from 2021-10-04 to 2021-10-05
$ws.Name = "Through 2021-10-05"

# Update the October row label to reflect the new "through" date
$ws.Range("A11").Value = "October (through 10-05)"

# New carjacking count for 2021 (column H), row 2 (January) ticks up by one
$ws.Range("H2").Value = 217

# October (row 11) gets one additional incident recorded for each year column
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 7
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 12
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 29
$ws.Range("H11").Value = 35

# Totals row (row 12) reflect the updated October figures
$ws.Range("B12").Value = 231
$ws.Range("C12").Value = 436
$ws.Range("D12").Value = 638
$ws.Range("E12").Value = 560
$ws.Range("F12").Value = 426
$ws.Range("G12").Value = 930
$ws.Range("H12").Value = 1284
